$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column E (Description) to better fit the new, longer description text.
$ws.Columns.Item(5).ColumnWidth = 61.6

# Add the new bug report row (row 6).
$ws.Cells.Item(6, 1).Value2 = 5
$ws.Cells.Item(6, 2).Value2 = "Images missing from PDF output"
$ws.Cells.Item(6, 3).Value2 = "OPEN"
$ws.Cells.Item(6, 4).Value2 = "DocumentViewer"
$ws.Cells.Item(6, 5).Value2 = "The wkhtmltopdf application seem to break images with a 'file://' url when converting to PDFs. This is likely a bug that I can't easily fix, so a workaround is needed."
$ws.Cells.Item(6, 6).Value2 = 40247
$ws.Cells.Item(6, 7).Value2 = 40247

# Match the row height used by similarly-wrapped rows (e.g. row 5).
$ws.Rows.Item(6).RowHeight = 45

# Move the selection down to the next empty row, like after data entry.
$ws.Range("A7").Select() | Out-Null
